$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 70, shifting existing rows 70-152 down to 71-153
$ws.Rows.Item(70).Insert()

# Populate the newly inserted row 70 with the new data record
$ws.Cells.Item(70, 1).Value = 11
$ws.Cells.Item(70, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(70, 3).Value = "Bíobío"
$ws.Cells.Item(70, 4).Value = 44638
$ws.Cells.Item(70, 5).Value = 8
$ws.Cells.Item(70, 6).Value = 100112003
$ws.Cells.Item(70, 7).Value = "Ajo"
$ws.Cells.Item(70, 8).Value = "Chino"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 220
$ws.Cells.Item(70, 11).Value = 17000
$ws.Cells.Item(70, 12).Value = 18000
$ws.Cells.Item(70, 13).Value = 17545
$ws.Cells.Item(70, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(70, 15).Value = "China"
$ws.Cells.Item(70, 16).Value = 1754
$ws.Cells.Item(70, 17).Value = 10
$ws.Cells.Item(70, 18).Value = "Hortaliza"
